$d = $word.ActiveDocument

$pairs = @(
    @("36×35=", "53×45="),
    @("55×39=", "57×93="),
    @("74×96=", "48×42="),
    @("36×62=", "36×83="),
    @("25×74=", "37×30="),
    @("20×36=", "91×75="),
    @("24×75=", "72×30="),
    @("36×53=", "63×96="),
    @("51×33=", "98×81="),
    @("35×11=", "18×56="),
    @("47×98=", "80×91="),
    @("44×71=", "48×76="),
    @("14×39=", "71×58="),
    @("60×39=", "69×99="),
    @("37×69=", "96×64="),
    @("85×12=", "58×91="),
    @("17×69=", "15×16="),
    @("81×92=", "17×76="),
    @("25×54=", "67×44="),
    @("94×68=", "46×71="),
    @("42×66=", "44×74="),
    @("56×71=", "21×87="),
    @("78×86=", "29×50="),
    @("82×72=", "17×36="),
    @("12×19=", "63×75=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
